$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# New regenerated values (B:E, G) for rows 2-7. F (Win) column is unchanged.
$ws.Range("B2").Value = 3.272327238179451
$ws.Range("C2").Value = 0.3048912486333797
$ws.Range("D2").Value = 0.1496068669990043
$ws.Range("E2").Value = 0.5333859586016987
$ws.Range("G2").Value = 4.260211312413533

$ws.Range("B3").Value = 3.272327238179451
$ws.Range("C3").Value = 1.626987699542094
$ws.Range("D3").Value = 0.7210945179870265
$ws.Range("E3").Value = 0.5333859586016987
$ws.Range("G3").Value = 6.15379541431027

$ws.Range("B4").Value = 3.272327238179451
$ws.Range("C4").Value = 0.3048912486333797
$ws.Range("D4").Value = 0.1496068669990043
$ws.Range("E4").Value = 0.5333859586016987
$ws.Range("G4").Value = 4.260211312413533

$ws.Range("B5").Value = 3.272327238179451
$ws.Range("C5").Value = 1.626987699542094
$ws.Range("D5").Value = 0.7210945179870265
$ws.Range("E5").Value = 0.5333859586016987
$ws.Range("G5").Value = 6.15379541431027

$ws.Range("B6").Value = 3.272327238179451
$ws.Range("C6").Value = 1.626987699542094
$ws.Range("D6").Value = 0.7210945179870265
$ws.Range("E6").Value = 0.5333859586016987
$ws.Range("G6").Value = 6.15379541431027

$ws.Range("B7").Value = 3.272327238179451
$ws.Range("C7").Value = 1.626987699542094
$ws.Range("D7").Value = 0.7210945179870265
$ws.Range("E7").Value = 0.5333859586016987
$ws.Range("G7").Value = 6.15379541431027
